# Add "Percentage" and "Rank" columns (AL, AM) to the NURSERY result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Copy the formatting of the existing "Total Marks" header (AK1, bold + border + centered)
# onto the two new header cells, then set their text.
$ws.Range("AK1").Copy()
$ws.Range("AL1:AM1").PasteSpecial(-4122)
$ws.Range("AL1").Value = "Percentage"
$ws.Range("AM1").Value = "Rank"

# --- Data rows (2-21) ---
# Column AI holds "Total Weightage" and AK holds the constant "Total Marks" (600).
# Percentage = Total Weightage / Total Marks * 100
# Rank = standard rank (1 = highest) of the Percentage values among all students,
# matching Excel's RANK(value, range, 0) semantics (ties share the same rank).

$firstRow = 2
$lastRow = 21

# First pass: read Total Weightage / Total Marks and compute each student's percentage.
$percentages = @{}
foreach ($r in $firstRow..$lastRow) {
    $totalWeightage = $ws.Range("AI$r").Value2
    $totalMarks = $ws.Range("AK$r").Value2
    $percentages[$r] = $totalWeightage / $totalMarks * 100
}

# Second pass: rank (1 = highest percentage) and write both new columns.
foreach ($r in $firstRow..$lastRow) {
    $myPct = $percentages[$r]
    $rank = 1
    foreach ($other in $firstRow..$lastRow) {
        if ($percentages[$other] -gt $myPct) {
            $rank = $rank + 1
        }
    }
    $ws.Range("AL$r").Value = $myPct
    $ws.Range("AM$r").Value = $rank
}
